# SenderRank.xlsx — "Add files via upload" re-upload.
# The uploaded version re-sorts/relabels the diplomatic-rank lookup table
# (column A label / column B count) starting at row 7, and appends one new
# row (37) for a straight-apostrophe "Chargé d'affaires a.i." variant that
# wasn't present before. Rows 1-6 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column A (label) / column B (count) pairs for rows 7..37, in final order.
$rows = @(
    @{ Row = 7;  Label = "papal nuncio";                   Count = 1 },
    @{ Row = 8;  Label = "Papal nuncio";                   Count = 1 },
    @{ Row = 9;  Label = "Papal Nuncio";                   Count = 1 },
    @{ Row = 10; Label = "pro-nuncio";                     Count = 1 },
    @{ Row = 11; Label = "Pro nuncio";                     Count = 1 },
    @{ Row = 12; Label = "Pro-Nuncio";                     Count = 1 },
    @{ Row = 13; Label = "Pro-nuncio";                     Count = 1 },
    @{ Row = 14; Label = "High commissioner";              Count = 1 },
    @{ Row = 15; Label = "High Commissioner";              Count = 1 },
    @{ Row = 16; Label = "Secretary of Peoples Bureau";    Count = 1 },
    @{ Row = 17; Label = "Minister";                       Count = 1 },
    @{ Row = 18; Label = "Envoy";                           Count = 2 },
    @{ Row = 19; Label = "Internuncio";                     Count = 2 },
    @{ Row = 20; Label = "Chargé d’affaires a.i";           Count = 2 },
    @{ Row = 21; Label = "Chargé d’affaires a.i.";          Count = 3 },
    @{ Row = 22; Label = "Chargé d'affaires a.i.";          Count = 3 },
    @{ Row = 23; Label = "Charge d’aﬀairs a.i.";            Count = 3 },
    @{ Row = 24; Label = "Charge d’aﬀairs ad hoc";          Count = 3 },
    @{ Row = 25; Label = "Charge d’aﬀairs ad interim";      Count = 3 },
    @{ Row = 26; Label = "Chargé d’affaires en titre";      Count = 3 },
    @{ Row = 27; Label = "Charge d’aﬀairs";                 Count = 4 },
    @{ Row = 28; Label = "Charge d’aﬀaires et pied";        Count = 4 },
    @{ Row = 29; Label = "Consul";                          Count = 5 },
    @{ Row = 30; Label = "Consul General";                  Count = 5 },
    @{ Row = 31; Label = "...";                             Count = 6 },
    @{ Row = 32; Label = "De facto diplomatic mission";     Count = 7 },
    @{ Row = 33; Label = "Others";                          Count = 8 },
    @{ Row = 34; Label = "Apostolic delegate";              Count = 8 },
    @{ Row = 35; Label = "commissioner";                    Count = 8 },
    @{ Row = 36; Label = "Commissioner";                    Count = 8 },
    @{ Row = 37; Label = "Unknown or missing";              Count = 9 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Label
    $ws.Cells.Item($r.Row, 2).Value = $r.Count
}

# Window/view state: zoom + new active selection (matches the uploaded file).
[void]$ws.Activate()
$excel.ActiveWindow.Zoom = 96
[void]$ws.Range("A20").Select()
